# Apply "contingencies with rene fine" edit:
# - Insert two new line rows (line7, line8) right after the existing line6 row,
#   pushing the extr1..extr8 rows down by two rows.
# - Update some of the numeric C/D/E values for the shifted rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-HeaderLikeStyle($rng) {
    # Mirrors the bold/centered/bordered style ("s=1") used by every cell in
    # column A of the data rows.
    $rng.Font.Bold = $true
    $rng.HorizontalAlignment = -4108   # xlCenter
    $rng.VerticalAlignment = -4160     # xlTop
    $rng.Borders.LineStyle = 1
}

# Final target contents for rows 2..17 (row, A, B, C, D, E)
$rows = @(
    @(2,  0,  "line1", 7,  9,  $true),
    @(3,  1,  "line2", 9,  8,  $true),
    @(4,  2,  "line3", 8,  10, $true),
    @(5,  3,  "line4", 8,  11, $true),
    @(6,  4,  "line5", 10, 5,  $false),
    @(7,  5,  "line6", 12, 8,  $true),
    @(8,  6,  "line7", 14, 11, $true),
    @(9,  7,  "line8", 16, 9,  $true),
    @(10, 8,  "extr1", 5,  12, $true),
    @(11, 9,  "extr2", 5,  9,  $true),
    @(12, 10, "extr3", 10, 11, $false),
    @(13, 11, "extr4", 7,  8,  $false),
    @(14, 12, "extr5", 9,  11, $false),
    @(15, 13, "extr6", 7,  11, $false),
    @(16, 14, "extr7", 5,  7,  $false),
    @(17, 15, "extr8", 8,  5,  $true)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
    $ws.Range("E$r").Value = $row[5]
    Set-HeaderLikeStyle($ws.Range("A$r"))
}
